# Applies the "Add files via upload" revision to CORE_holdings.xlsx:
#  1) Bumps the "as of" date in the confidentiality disclosure from
#     2021-03-30 to 2021-03-31.
#  2) Refreshes the Weight (column D) and Percent Change (column E)
#     figures for rows 2-8 to their newly-recalculated values.
#
# The sheet is protected, so it has to be unprotected before editing and
# re-protected (with the same password) afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# --- Disclosure text: bump the "as of" date by one day -------------------
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."
# Re-run autofit so the row keeps its original (default) height instead of
# picking up an explicit custom height from the multi-line text re-entry.
$ws.Rows.Item(11).AutoFit()

# --- Weight / Percent Change refresh (rows 2-8) ---------------------------
$ws.Range("D2").Value = 0.4995001934638222
$ws.Range("E2").Value = -0.004861551469033998

$ws.Range("D3").Value = 0.2418355767471482
$ws.Range("E3").Value = 0.0121249805689414

$ws.Range("D4").Value = 0.09803486747978127
$ws.Range("E4").Value = 0.007067591878694435

$ws.Range("D5").Value = 0.1025008323443516
$ws.Range("E5").Value = -0.005754974639094801

$ws.Range("D6").Value = 0.03025534225099906
$ws.Range("E6").Value = -0.001288787548329728

$ws.Range("D7").Value = 0.0278731877138975
$ws.Range("E7").Value = 0.01191603535353525

$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = 0.0009000316940730446

$ws.Protect("D382")
